$d = $word.ActiveDocument

# Locate the paragraph that ends with the dev-branch sentence.
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*创建了一个dev分支*") {
        $target = $p
        break
    }
}

if ($target -ne $null) {
    $r = $target.Range
    # Move to just before the paragraph mark (end of the text, before the pilcrow).
    $insertRange = $r.Duplicate
    $insertRange.Start = $r.End - 1
    $insertRange.End = $r.End - 1
    $insertRange.Text = "使用git创建分支简单快捷。"
    $insertRange.Font.NameFarEast = "宋体"
}
